# Merge split runs back into single runs (COM Find/Replace naturally
# collapses a multi-run match into one run carrying the replacement text).

$d = $word.ActiveDocument

function Merge-Text($text) {
    $rng = $d.Content
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

# "- Project 3 is due June 21st so we have 1 " + "set " + "meeting after this to get it done"
Merge-Text "- Project 3 is due June 21st so we have 1 set meeting after this to get it done"

# "Fixing the Branch and truncation problems from milestone " + "2"  (two occurrences)
Merge-Text "Fixing the Branch and truncation problems from milestone 2"

# " 6/1" + "2" + "/2024"
Merge-Text " 6/12/2024"

# "Class Design rough draft" + " (All except the UI)"
Merge-Text "Class Design rough draft (All except the UI)"

# "Class Design rough draft (" + "only " + "UI)"
Merge-Text "Class Design rough draft (only UI)"

# "Simplify " + "GUI Design"
Merge-Text "Simplify GUI Design"

# "Work on team SRS doc" + ", Complete before Tuesday meeting" + " " + "("
Merge-Text "Work on team SRS doc, Complete before Tuesday meeting ("

# Append the new second-meeting report (6/18/2024) after the last paragraph.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$newXml = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>DATE:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve"> 6/1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>8</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>/2024</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>ATTENDANCE:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="321242054"/>
          <w14:checkbox>
            <w14:checked w14:val="1"/>
            <w14:checkedState w14:val="2612" w14:font="MS Gothic"/>
            <w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/>
          </w14:checkbox>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/>
            </w:rPr>
            <w:t>☒</w:t>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t xml:space="preserve">   Andrew Buckland</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="-342938904"/>
          <w14:checkbox>
            <w14:checked w14:val="1"/>
            <w14:checkedState w14:val="2612" w14:font="MS Gothic"/>
            <w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/>
          </w14:checkbox>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/>
            </w:rPr>
            <w:t>☒</w:t>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t xml:space="preserve">   </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Jaykant</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Chaudhary Lekhi</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="1613623958"/>
          <w14:checkbox>
            <w14:checked w14:val="1"/>
            <w14:checkedState w14:val="2612" w14:font="MS Gothic"/>
            <w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/>
          </w14:checkbox>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/>
            </w:rPr>
            <w:t>☒</w:t>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t xml:space="preserve">   Jonah Cragun</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="-1009678346"/>
          <w14:checkbox>
            <w14:checked w14:val="1"/>
            <w14:checkedState w14:val="2612" w14:font="MS Gothic"/>
            <w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/>
          </w14:checkbox>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/>
            </w:rPr>
            <w:t>☒</w:t>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t xml:space="preserve">   Zac Peterson</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>AGENDA:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Recap on current action items</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Plan for final push</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>ACTION ITEMS:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>ANDREW</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>JAYKANT</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>JONAH</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Look into possibilities of using </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>web based</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> GUI instead of Qt</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>ZAC</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Figure out how to implement the code into the Qt GUI and how to build</w:t>
      </w:r>
    </w:p>
</w:body>
'@
$endRange.InsertXML($newXml)
